$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.354.54"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "2.066.02"
$ws.Range("E3").Value = "  +3.32%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5171"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4290"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.146"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").Value = "2.067.28"
$ws.Range("E13").Value = "  +3.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.581"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.606"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001107"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06602"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.185"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "30.422.38"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("D26").Value = "2.317.15"
$ws.Range("E26").Value = "  +3.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.483"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.174"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1062"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.006"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.808"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.482"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02540"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.485"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.420"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06560"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2224"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6638"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.229"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6221"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.176"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("E49").Value = "  -3.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.174"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.68%  "
